# Weekly update: insert a new sampling row for "Macroferia Regional de
# Talca - Zanahoria" just before the current row 324, shifting the
# existing rows 324:345 down to 325:346 (the sheet's dimension grows
# from A1:R345 to A1:R346).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 324..345 down one position, opening up a blank row 324.
$ws.Rows.Item(324).Insert()

# Populate the newly-inserted row 324 with this week's observation.
$ws.Cells.Item(324, 1).Value  = 5
$ws.Cells.Item(324, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(324, 3).Value  = "Maule"
$ws.Cells.Item(324, 4).Value  = 44714
$ws.Cells.Item(324, 5).Value  = 7
$ws.Cells.Item(324, 6).Value  = 100114013
$ws.Cells.Item(324, 7).Value  = "Zanahoria"
$ws.Cells.Item(324, 8).Value  = "Sin especificar"
$ws.Cells.Item(324, 9).Value  = "Primera"
$ws.Cells.Item(324, 10).Value = 500
$ws.Cells.Item(324, 11).Value = 5500
$ws.Cells.Item(324, 12).Value = 5500
$ws.Cells.Item(324, 13).Value = 5500
$ws.Cells.Item(324, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(324, 15).Value = "Región de Ñuble"
$ws.Cells.Item(324, 16).Value = 275
$ws.Cells.Item(324, 17).Value = 20
$ws.Cells.Item(324, 18).Value = "Hortaliza"
